$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the item row with the new item's data ---
$ws.Range("A2").Value2 = "Souvenir SSG 08 | Jungle Dashed (Battle-Scarred)"
$ws.Range("B2").Value2 = "def;0.45-0.46"
$ws.Range("C2").Value2 = "0-0.8;0-0.1"
$ws.Range("D2").Value2 = 1

# --- Row 3: the second item row is removed, leaving only the styled blank B3 ---
$ws.Range("A3:D3").ClearContents()

# --- Remove the 4 trailing blank rows (48:51) ---
$ws.Rows("48:51").Delete()

# --- Column widths / new column D ---
# (Target character-widths are 47.140625 / 22.140625 / 15 / 8.140625; the host
# rounds ColumnWidth to the nearest 1/6 of a character, so we feed it
# "target - 5/6" which lands on the closest representable step.)
$ws.Columns("A").ColumnWidth = 46.307291666666664
$ws.Columns("B").ColumnWidth = 21.307291666666668
$ws.Columns("C").ColumnWidth = 14.166666666666666
$ws.Columns("D").ColumnWidth = 7.307291666666667

# --- Selection moves to D5 ---
$ws.Range("D5").Select()

Write-Output "edit applied"
